$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.460.91'
$ws.Range('E2').Value = '  -1.81%  '
$ws.Range('D3').Value = '2.509.18'
$ws.Range('E3').Value = '  -4.91%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.96'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '170.80'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.525'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = '2.508.33'
$ws.Range('E9').Value = '  -4.90%  '
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.12'
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.349'
$ws.Range('E13').Value = '  -4.73%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.67'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('D15').Value = '2.957.25'
$ws.Range('E16').Value = '  -3.55%  '
$ws.Range('D17').Value = '66.238.82'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '2.490.11'
$ws.Range('E18').Value = '  -6.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.25'
$ws.Range('E19').Value = '  -7.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.71'
$ws.Range('E20').Value = '  -4.70%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '347.51'
$ws.Range('E21').Value = '  -3.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.19'
$ws.Range('E22').Value = '  -3.50%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.63'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '69.80'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('E27').Value = '  -4.70%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('E29').Value = '  -5.05%  '
$ws.Range('D30').Value = '0.0₃0975'
$ws.Range('E30').Value = '  -3.85%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '523.49'
$ws.Range('E31').Value = '  -4.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.07'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.32'
$ws.Range('E33').Value = '  -3.51%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  -3.74%  '
$ws.Range('E35').Value = '  -4.94%  '
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '157.02'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.62'
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.38'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.355'
$ws.Range('E41').Value = '  -3.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.79'
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.08'
$ws.Range('E43').Value = '  -3.27%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '39.33'
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '148.56'
$ws.Range('E47').Value = '  -3.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.557'
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.68'
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  -10.92%  '
